$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 35.88321066666667
$ws.Range("H2").Value = 107.649632
$ws.Range("I2").Value = 0.08317795499144418
$ws.Range("J2").Value = 0.08448843719082051
$ws.Range("M2").Value = 13.89934866666667
$ws.Range("N2").Value = 41.69804600000001
$ws.Range("O2").Value = 0.04853507553134179
$ws.Range("P2").Value = 0.04999273878390351
$ws.Range("Q2").Value = 498.7532563354525
$ws.Range("R2").Value = 4488.779307019073
$ws.Range("S2").Value = 0.004037048328052291
$ws.Range("T2").Value = 0.004223808370740928
# Row 3
$ws.Range("G3").Value = 35.88321066666667
$ws.Range("H3").Value = 107.649632
$ws.Range("I3").Value = 0.08317795499144418
$ws.Range("J3").Value = 0.08448843719082051
$ws.Range("O3").Value = 0.245697991654417
$ws.Range("P3").Value = 0.253077086664408
$ws.Range("Q3").Value = 2524.827087857085
$ws.Range("R3").Value = 22723.44379071376
$ws.Range("S3").Value = 0.02043665649131933
$ws.Range("T3").Value = 0.02138208754108167
# Row 4
$ws.Range("G4").Value = 35.88321066666667
$ws.Range("H4").Value = 107.649632
$ws.Range("I4").Value = 0.08317795499144418
$ws.Range("J4").Value = 0.08448843719082051
$ws.Range("M4").Value = 82.007665
$ws.Range("N4").Value = 246.022995
$ws.Range("O4").Value = 0.2863622109480123
$ws.Range("P4").Value = 0.2949625822722868
$ws.Range("Q4").Value = 2942.698319476427
$ws.Range("R4").Value = 26484.28487528785
$ws.Range("S4").Value = 0.02381902309348421
$ws.Range("T4").Value = 0.02492092760595433
# Row 5
$ws.Range("G5").Value = 35.88321066666667
$ws.Range("H5").Value = 107.649632
$ws.Range("I5").Value = 0.08317795499144418
$ws.Range("J5").Value = 0.08448843719082051
$ws.Range("M5").Value = 25.0501465
$ws.Range("N5").Value = 50.100293
$ws.Range("O5").Value = 0.0874724982879541
$ws.Range("P5").Value = 0.06006638442832619
$ws.Range("Q5").Value = 898.8796840903627
$ws.Range("R5").Value = 5393.278104542177
$ws.Range("S5").Value = 0.007275783525584624
$ws.Range("T5").Value = 0.005074914948052316
# Row 6
$ws.Range("G6").Value = 35.88321066666667
$ws.Range("H6").Value = 107.649632
$ws.Range("I6").Value = 0.08317795499144418
$ws.Range("J6").Value = 0.08448843719082051
$ws.Range("M6").Value = 95.05788666666668
$ws.Range("N6").Value = 285.17366
$ws.Range("O6").Value = 0.3319322235782747
$ws.Range("P6").Value = 0.3419012078510756
$ws.Range("Q6").Value = 3410.982172788125
$ws.Range("R6").Value = 30698.83955509313
$ws.Range("S6").Value = 0.02760944355300372
$ws.Range("T6").Value = 0.02888669872499127
# Row 7
$ws.Range("I7").Value = 0.03522729558434242
$ws.Range("J7").Value = 0.03578230735158529
$ws.Range("M7").Value = 13.89934866666667
$ws.Range("N7").Value = 41.69804600000001
$ws.Range("O7").Value = 0.04853507553134179
$ws.Range("P7").Value = 0.04999273878390351
$ws.Range("Q7").Value = 211.2305885181902
$ws.Range("R7").Value = 1901.075296663712
$ws.Range("S7").Value = 0.001709759451950962
$ws.Range("T7").Value = 0.001788855544513154
# Row 8
$ws.Range("I8").Value = 0.03522729558434242
$ws.Range("J8").Value = 0.03578230735158529
$ws.Range("O8").Value = 0.245697991654417
$ws.Range("P8").Value = 0.253077086664408
$ws.Range("S8").Value = 0.008655275776489446
$ws.Range("T8").Value = 0.009055682098669634
# Row 9
$ws.Range("I9").Value = 0.03522729558434242
$ws.Range("J9").Value = 0.03578230735158529
$ws.Range("M9").Value = 82.007665
$ws.Range("N9").Value = 246.022995
$ws.Range("O9").Value = 0.2863622109480123
$ws.Range("P9").Value = 0.2949625822722868
$ws.Range("Q9").Value = 1246.283387544293
$ws.Range("R9").Value = 11216.55048789864
$ws.Range("S9").Value = 0.01008776624925145
$ws.Range("T9").Value = 0.01055444177608423
# Row 10
$ws.Range("I10").Value = 0.03522729558434242
$ws.Range("J10").Value = 0.03578230735158529
$ws.Range("M10").Value = 25.0501465
$ws.Range("N10").Value = 50.100293
$ws.Range("O10").Value = 0.0874724982879541
$ws.Range("P10").Value = 0.06006638442832619
$ws.Range("Q10").Value = 380.6910175835493
$ws.Range("R10").Value = 2284.146105501296
$ws.Range("S10").Value = 0.003081419552690646
$ws.Range("T10").Value = 0.002149313829112845
# Row 11
$ws.Range("I11").Value = 0.03522729558434242
$ws.Range("J11").Value = 0.03578230735158529
$ws.Range("M11").Value = 95.05788666666668
$ws.Range("N11").Value = 285.17366
$ws.Range("O11").Value = 0.3319322235782747
$ws.Range("P11").Value = 0.3419012078510756
$ws.Range("Q11").Value = 1444.609659447502
$ws.Range("R11").Value = 13001.48693502752
$ws.Range("S11").Value = 0.01169307455395992
$ws.Range("T11").Value = 0.01223401410320544
# Row 12
$ws.Range("G12").Value = 177.70077
$ws.Range("H12").Value = 533.10231
$ws.Range("I12").Value = 0.4119137160358794
$ws.Range("J12").Value = 0.4184034835782469
$ws.Range("M12").Value = 13.89934866666667
$ws.Range("N12").Value = 41.69804600000001
$ws.Range("O12").Value = 0.04853507553134179
$ws.Range("P12").Value = 0.04999273878390351
$ws.Range("Q12").Value = 2469.92496056514
$ws.Range("R12").Value = 22229.32464508626
$ws.Range("S12").Value = 0.01999226332019708
$ws.Range("T12").Value = 0.02091713606080256
# Row 13
$ws.Range("G13").Value = 177.70077
$ws.Range("H13").Value = 533.10231
$ws.Range("I13").Value = 0.4119137160358794
$ws.Range("J13").Value = 0.4184034835782469
$ws.Range("O13").Value = 0.245697991654417
$ws.Range("P13").Value = 0.253077086664408
$ws.Range("Q13").Value = 12503.44407017745
$ws.Range("R13").Value = 112530.996631597
$ws.Range("S13").Value = 0.1012063727649234
$ws.Range("T13").Value = 0.1058883346742222
# Row 14
$ws.Range("G14").Value = 177.70077
$ws.Range("H14").Value = 533.10231
$ws.Range("I14").Value = 0.4119137160358794
$ws.Range("J14").Value = 0.4184034835782469
$ws.Range("M14").Value = 82.007665
$ws.Range("N14").Value = 246.022995
$ws.Range("O14").Value = 0.2863622109480123
$ws.Range("P14").Value = 0.2949625822722868
$ws.Range("Q14").Value = 14572.82521640205
$ws.Range("R14").Value = 131155.4269476185
$ws.Range("S14").Value = 0.1179565224438461
$ws.Range("T14").Value = 0.12341337194796
# Row 15
$ws.Range("G15").Value = 177.70077
$ws.Range("H15").Value = 533.10231
$ws.Range("I15").Value = 0.4119137160358794
$ws.Range("J15").Value = 0.4184034835782469
$ws.Range("M15").Value = 25.0501465
$ws.Range("N15").Value = 50.100293
$ws.Range("O15").Value = 0.0874724982879541
$ws.Range("P15").Value = 0.06006638442832619
$ws.Range("Q15").Value = 4451.430321662805
$ws.Range("R15").Value = 26708.58192997683
$ws.Range("S15").Value = 0.03603112182073327
$ws.Range("T15").Value = 0.02513198449076184
# Row 16
$ws.Range("G16").Value = 177.70077
$ws.Range("H16").Value = 533.10231
$ws.Range("I16").Value = 0.4119137160358794
$ws.Range("J16").Value = 0.4184034835782469
$ws.Range("M16").Value = 95.05788666666668
$ws.Range("N16").Value = 285.17366
$ws.Range("O16").Value = 0.3319322235782747
$ws.Range("P16").Value = 0.3419012078510756
$ws.Range("Q16").Value = 16891.8596552394
$ws.Range("R16").Value = 152026.7368971546
$ws.Range("S16").Value = 0.1367274356861795
$ws.Range("T16").Value = 0.1430526564045003
# Row 17
$ws.Range("G17").Value = 20.074196
$ws.Range("H17").Value = 40.148392
$ws.Range("I17").Value = 0.04653236263856699
$ws.Range("J17").Value = 0.0315103250497358
$ws.Range("M17").Value = 13.89934866666667
$ws.Range("N17").Value = 41.69804600000001
$ws.Range("O17").Value = 0.04853507553134179
$ws.Range("P17").Value = 0.04999273878390351
$ws.Range("Q17").Value = 279.0182494070054
$ws.Range("R17").Value = 1674.109496442032
$ws.Range("S17").Value = 0.002258451735314635
$ws.Range("T17").Value = 0.001575287449207333
# Row 18
$ws.Range("G18").Value = 20.074196
$ws.Range("H18").Value = 40.148392
$ws.Range("I18").Value = 0.04653236263856699
$ws.Range("J18").Value = 0.0315103250497358
$ws.Range("O18").Value = 0.245697991654417
$ws.Range("P18").Value = 0.253077086664408
$ws.Range("Q18").Value = 1412.467638377593
$ws.Range("R18").Value = 8474.805830265559
$ws.Range("S18").Value = 0.01143290804723094
$ws.Range("T18").Value = 0.007974541263435652
# Row 19
$ws.Range("G19").Value = 20.074196
$ws.Range("H19").Value = 40.148392
$ws.Range("I19").Value = 0.04653236263856699
$ws.Range("J19").Value = 0.0315103250497358
$ws.Range("M19").Value = 82.007665
$ws.Range("N19").Value = 246.022995
$ws.Range("O19").Value = 0.2863622109480123
$ws.Range("P19").Value = 0.2949625822722868
$ws.Range("Q19").Value = 1646.23794071234
$ws.Range("R19").Value = 9877.427644274041
$ws.Range("S19").Value = 0.01332511024581473
$ws.Range("T19").Value = 0.009294366844909195
# Row 20
$ws.Range("G20").Value = 20.074196
$ws.Range("H20").Value = 40.148392
$ws.Range("I20").Value = 0.04653236263856699
$ws.Range("J20").Value = 0.0315103250497358
$ws.Range("M20").Value = 25.0501465
$ws.Range("N20").Value = 50.100293
$ws.Range("O20").Value = 0.0874724982879541
$ws.Range("P20").Value = 0.06006638442832619
$ws.Range("Q20").Value = 502.861550669714
$ws.Range("R20").Value = 2011.446202678856
$ws.Range("S20").Value = 0.00407030201123651
$ws.Range("T20").Value = 0.001892711297898947
# Row 21
$ws.Range("G21").Value = 20.074196
$ws.Range("H21").Value = 40.148392
$ws.Range("I21").Value = 0.04653236263856699
$ws.Range("J21").Value = 0.0315103250497358
$ws.Range("M21").Value = 95.05788666666668
$ws.Range("N21").Value = 285.17366
$ws.Range("O21").Value = 0.3319322235782747
$ws.Range("P21").Value = 0.3419012078510756
$ws.Range("Q21").Value = 1908.210648292453
$ws.Range("R21").Value = 11449.26388975472
$ws.Range("S21").Value = 0.01544559059897017
$ws.Range("T21").Value = 0.01077341819428468
# Row 22
$ws.Range("G22").Value = 182.547562
$ws.Range("H22").Value = 547.642686
$ws.Range("I22").Value = 0.423148670749767
$ws.Range("J22").Value = 0.4298154468296114
$ws.Range("M22").Value = 13.89934866666667
$ws.Range("N22").Value = 41.69804600000001
$ws.Range("O22").Value = 0.04853507553134179
$ws.Range("P22").Value = 0.04999273878390351
$ws.Range("Q22").Value = 2537.292212487951
$ws.Range("R22").Value = 22835.62991239156
$ws.Range("S22").Value = 0.02053755269582682
$ws.Range("T22").Value = 0.02148765135863953
# Row 23
$ws.Range("G23").Value = 182.547562
$ws.Range("H23").Value = 547.642686
$ws.Range("I23").Value = 0.423148670749767
$ws.Range("J23").Value = 0.4298154468296114
$ws.Range("O23").Value = 0.245697991654417
$ws.Range("P23").Value = 0.253077086664408
$ws.Range("Q23").Value = 12844.47575333664
$ws.Range("R23").Value = 115600.2817800297
$ws.Range("S23").Value = 0.1039667785744539
$ws.Range("T23").Value = 0.1087764410869988
# Row 24
$ws.Range("G24").Value = 182.547562
$ws.Range("H24").Value = 547.642686
$ws.Range("I24").Value = 0.423148670749767
$ws.Range("J24").Value = 0.4298154468296114
$ws.Range("M24").Value = 82.007665
$ws.Range("N24").Value = 246.022995
$ws.Range("O24").Value = 0.2863622109480123
$ws.Range("P24").Value = 0.2949625822722868
$ws.Range("Q24").Value = 14970.29931106273
$ws.Range("R24").Value = 134732.6937995646
$ws.Range("S24").Value = 0.1211737889156158
$ws.Range("T24").Value = 0.126779474097379
# Row 25
$ws.Range("G25").Value = 182.547562
$ws.Range("H25").Value = 547.642686
$ws.Range("I25").Value = 0.423148670749767
$ws.Range("J25").Value = 0.4298154468296114
$ws.Range("M25").Value = 25.0501465
$ws.Range("N25").Value = 50.100293
$ws.Range("O25").Value = 0.0874724982879541
$ws.Range("P25").Value = 0.06006638442832619
$ws.Range("Q25").Value = 4572.843171317833
$ws.Range("R25").Value = 27437.059027907
$ws.Range("S25").Value = 0.03701387137770905
$ws.Range("T25").Value = 0.02581745986250024
# Row 26
$ws.Range("G26").Value = 182.547562
$ws.Range("H26").Value = 547.642686
$ws.Range("I26").Value = 0.423148670749767
$ws.Range("J26").Value = 0.4298154468296114
$ws.Range("M26").Value = 95.05788666666668
$ws.Range("N26").Value = 285.17366
$ws.Range("O26").Value = 0.3319322235782747
$ws.Range("P26").Value = 0.3419012078510756
$ws.Range("Q26").Value = 17352.58545987231
$ws.Range("R26").Value = 156173.2691388508
$ws.Range("S26").Value = 0.1404566791861614
